$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.676.25"
$ws.Range("E2").Value = "  +2.26%  "

# Row 3
$ws.Range("D3").Value = "3.554.10"
$ws.Range("E3").Value = "  +1.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'580.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "

# Row 6
$ws.Range("D6").Value = "'186.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.03%  "

# Row 7
$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.81%  "

# Row 8
$ws.Range("D8").Value = "3.546.33"
$ws.Range("E8").Value = "  +1.29%  "

# Row 9
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("E10").Value = "  +18.71%  "

# Row 11
$ws.Range("D11").Value = "'0.648"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.66%  "

# Row 12
$ws.Range("D12").Value = "'54.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "

# Row 13
$ws.Range("E13").Value = "  +5.50%  "

# Row 14
$ws.Range("D14").Value = "'9.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.25%  "

# Row 15
$ws.Range("D15").Value = "4.128.74"
$ws.Range("E15").Value = "  +1.47%  "

# Row 16
$ws.Range("D16").Value = "70.768.87"
$ws.Range("E16").Value = "  +2.58%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'19.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.47%  "

# Row 18
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.52%  "

# Row 19
$ws.Range("D19").Value = "3.557.25"
$ws.Range("E19").Value = "  +1.42%  "

# Row 20
$ws.Range("D20").Value = "'573.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.78%  "

# Row 22
$ws.Range("E22").Value = "  -1.31%  "

# Row 23
$ws.Range("D23").Value = "'17.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.17%  "

# Row 24
$ws.Range("D24").Value = "'4.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.93%  "

# Row 25
$ws.Range("D25").Value = "'4.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "

# Row 26
$ws.Range("D26").Value = "'94.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.16%  "

# Row 27
$ws.Range("D27").Value = "'11.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.78%  "

# Row 28
$ws.Range("D28").Value = "'2.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.33%  "

# Row 29
$ws.Range("D29").Value = "'9.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.90%  "

# Row 30
$ws.Range("D30").Value = "'32.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.71%  "

# Row 31
$ws.Range("D31").Value = "'7.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.09%  "

# Row 32
$ws.Range("D32").Value = "'12.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.05%  "

# Row 33
$ws.Range("E33").Value = "  +1.61%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'63.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.11%  "

# Row 35
$ws.Range("B35").Value = "dogwifhat"
$ws.Range("C35").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D35").Value = "'3.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +21.39%  "

# Row 36
$ws.Range("E36").Value = "  +6.99%  "

# Row 37
$ws.Range("D37").Value = "'533.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.34%  "

# Row 38
$ws.Range("D38").Value = "'0.410"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.55%  "

# Row 39
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "'38.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.32%  "

# Row 40
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0806"
$ws.Range("E40").Value = "  +5.38%  "

# Row 42
$ws.Range("D42").Value = "3.624.49"
$ws.Range("E42").Value = "  +10.05%  "

# Row 43
$ws.Range("D43").Value = "'0.139"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.25%  "

# Row 44
$ws.Range("D44").Value = "'3.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.22%  "

# Row 45
$ws.Range("D45").Value = "'0.0468"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.34%  "

# Row 46
$ws.Range("E46").Value = "  +0.22%  "

# Row 47
$ws.Range("D47").Value = "'2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.91%  "

# Row 48
$ws.Range("D48").Value = "'9.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.92%  "

# Row 49
$ws.Range("E49").Value = "  +2.49%  "

# Row 50
$ws.Range("E50").Value = "  +0.18%  "

# Row 51
$ws.Range("D51").Value = "'1.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.91%  "
